# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.546.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.982.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.19%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '382.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.546'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0847'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.451.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.980.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.999'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.492.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0963'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '261.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.118'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +16.02%  '
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '34.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  -2.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0452'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.39%  '
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.91'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.274'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.032.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("E51").Value = '  +3.68%  '
